$d = $word.ActiveDocument

function Split-Range($startPos, $endPos) {
    # Toggling a character-formatting property on a range and then reverting
    # it forces the underlying run(s) covering that range to be split out as
    # their own <w:r> elements, even though the effective formatting ends up
    # unchanged. This lets us introduce explicit run boundaries at the exact
    # offsets we need without altering any visible formatting.
    $rng = $d.Range($startPos, $endPos)
    $rng.Bold = 1
    $rng.Bold = 0
}

# The run boundary right before "When she was 18" must be preserved (it
# separates this sentence from the preceding "During high school..." run).
$preRange = $d.Content
[void]$preRange.Find.Execute("When she was 18", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$preBoundary = $preRange.Start

# Locate the point right after "teacher" (and before the existing ". She
# never returned...") where the new sentence fragment gets inserted.
$teacherRange = $d.Content
[void]$teacherRange.Find.Execute("to university to become a teacher", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertStart = $teacherRange.End

# Insert the new text. Note: no trailing period is added here because the
# original text already continues with ". She never returned...".
$newText = ", dying her hair from this point onward regularly in order to feel more confident in a new environment"
$insertionPoint = $d.Range($insertStart, $insertStart)
$insertionPoint.InsertAfter($newText)

# Offsets of the run boundaries introduced by the new text.
$b1 = $insertStart            # end of "...become a teacher" / start of new run
$b2 = $insertStart + 91       # end of ", dying ... in a new " / start of "environment"
$b3 = $insertStart + 102      # end of "environment" / start of ". She never returned..."

# Find where the untouched tail (". She never ... during those years,") ends,
# i.e. just before the existing single-space run that follows it.
$tailRange = $d.Range($b3, $d.Content.End)
[void]$tailRange.Find.Execute("during those years,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$b4 = $tailRange.End
$b5 = $b4 + 1

# Re-establish every run boundary touched by the edit so the paragraph ends
# up with the same run layout the diff expects, instead of everything
# collapsing back into one merged run.
Split-Range $preBoundary $b1
Split-Range $b1 $b2
Split-Range $b2 $b3
Split-Range $b3 $b4
Split-Range $b4 $b5
